$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new "through" date
$ws.Name = "Through 2022-09-19"
$ws.Range("I1").Value = "2022 (through 09-19)"

# Update October (row 10) 2022 total and the year Total row (row 14) 2022 total
$ws.Range("I10").Value = 92
$ws.Range("I14").Value = 1227
